$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly time log entry: 10/27/2023, Internship, "Completed 8 hours assisting with daily operations" bucket
# Copy formatting from the row above (row 16) into the new row 17, then fill in values
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(17, 1).Value2 = 45226
$ws.Cells.Item(17, 2).Value2 = "Internship"
$ws.Cells.Item(17, 3).Value2 = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Mirror Excel's post-entry active cell landing on the next empty row
$ws.Range("C18").Select()
